$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "C2"  = "嘉美包装"
    "A3"  = "三花智控"
    "B3"  = "百达精工"
    "C3"  = "天奇股份"
    "A4"  = "嘉美包装"
    "C4"  = "博纳影业"
    "B5"  = "利亚德"
    "C5"  = "华胜天成"
    "A6"  = "博纳影业"
    "B6"  = "嘉美包装"
    "A7"  = "卧龙电驱"
    "B7"  = "三花智控"
    "C7"  = "风语筑"
    "A8"  = "华胜天成"
    "B8"  = "贵州茅台"
    "C8"  = "百达精工"
    "A9"  = "东方财富"
    "B9"  = "五洲新春"
    "C9"  = "汉缆股份"
    "A10" = "掌阅科技"
    "B10" = "中大力德"
    "C10" = "万向钱潮"
    "A11" = "光线传媒"
    "B11" = "东方财富"
    "C11" = "利欧股份"
    "A12" = "五洲新春"
    "B12" = "绿的谐波"
    "C12" = "协鑫集成"
    "B13" = "华胜天成"
    "C13" = "三花智控"
    "A14" = "风语筑"
    "B14" = "博纳影业"
    "C14" = "卧龙电驱"
    "A15" = "利欧股份"
    "B15" = "光线传媒"
    "C15" = "克来机电"
    "B16" = "深科技"
    "C16" = "大位科技"
    "A17" = "深科技"
    "B17" = "万向钱潮"
    "C17" = "巨力索具"
    "A18" = "汉缆股份"
    "B18" = "利欧股份"
    "C18" = "掌阅科技"
    "A19" = "贵州茅台"
    "B19" = "风语筑"
    "C19" = "深科技"
    "A20" = "中大力德"
    "B20" = "鸣志电器"
    "C20" = "五洲新春"
    "A21" = "通富微电"
    "B21" = "科大讯飞"
    "C21" = "国安股份"
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
